# Fix double entries for 2023-05-16 (date serial 45062) in the BL sheet.
# Rows 478-511 originally held the 2023-05-16 data duplicated (each of the
# 17 Bundesland rows appears twice). We delete the second occurrence of
# each duplicate pair; deleting entire rows shifts the following
# 2023-05-17 (45063) block up automatically, closing the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BL")

# Duplicate rows to remove: 479, 481, 483, ..., 511 (17 rows)
$rowsToDelete = 479..511 | Where-Object { ($_ % 2) -eq 1 }

# Delete from the bottom up so row indices above stay valid as we go.
$sorted = $rowsToDelete | Sort-Object -Descending
foreach ($r in $sorted) {
    $ws.Rows.Item($r).Delete()
}
